# Commit: "Created figure using NY sample"
#
# The "distance" sheet (the original first sheet, containing a small
# distance/HS/SL lookup table) is removed entirely. The remaining
# "Sheet1" (originally the second sheet) becomes the workbook's only
# sheet, and its header row is relabeled to the short column names used
# for the figure:
#   Location Latitude  -> lat
#   Location Longitude -> long
#   "Is the Hollywood sign illuminated/lit at night on a typical night?" -> holly_sign_acc
#   "Is she holding the torch in her left or right hand?"                -> statof_lib_acc
# (" ID" and "Zipcode" headers are unchanged.)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the now-unused "distance" sheet.
$distanceSheet = $wb.Worksheets.Item("distance")
$distanceSheet.Delete()

# Make "Sheet1" (the survivor) the active sheet.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Relabel the header cells for the accuracy / lat-long columns.
# (Set E1/F1 before B1/C1 so new shared-string entries are appended in
# the same order as the reference edit.)
$ws1.Range("E1").Value = "holly_sign_acc"
$ws1.Range("F1").Value = "statof_lib_acc"
$ws1.Range("B1").Value = "lat"
$ws1.Range("C1").Value = "long"

# Restore the sheet's prior selection on the new active cell.
$ws1.Range("H6").Select()

$excel.DisplayAlerts = $true
